# 03.07.2025 - maç sonuçları eklendi
# Add newly played match results (scores + highlight link) to the "Maçlar" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (Araklı 1961 Spor - Of 1461) now has its highlight video link.
$ws.Range("H16").Value = "https://youtu.be/LJGmM-PJz-g"

# Row 18 (Of FK - Ofside) final score: 17-3.
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 3

# Row 19 (Armedospor - Hubuş FK) final score: 3-7.
$ws.Range("F19").Value = 3
$ws.Range("G19").Value = 7

# Move the active selection to the next still-unplayed match (H17).
[void]$ws.Range("H17").Select()
